$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 399.25
$ws.Range("I4").Value = 399.25
$ws.Range("K4").Value = 399.25
$ws.Range("M4").Value = -285.25
$ws.Range("H29").Value = 1299.875
$ws.Range("I29").Value = 99.75
$ws.Range("J29").Value = 2500
$ws.Range("K29").Value = 299.25
$ws.Range("L29").Value = 7500
$ws.Range("M29").Value = -18.25
$ws.Range("N29").Value = -8062
$ws.Range("H31").Value = 1938.25
$ws.Range("I31").Value = 1938.25
$ws.Range("K31").Value = 5814.75
$ws.Range("M31").Value = -5584.75
$ws.Range("H32").Value = 5125
$ws.Range("I32").Value = 6500
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 6500
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -6174
$ws.Range("N32").Value = -1652
$ws.Range("H33").Value = 52696090
$ws.Range("I33").Value = 76924640
$ws.Range("J33").Value = 200895.5
$ws.Range("K33").Value = 76924640
$ws.Range("L33").Value = 200895.5
$ws.Range("M33").Value = -76924411
$ws.Range("N33").Value = -201353.5
$ws.Range("H34").Value = 25140
$ws.Range("I34").Value = 3350
$ws.Range("J34").Value = 39666.668
$ws.Range("K34").Value = 3350
$ws.Range("L34").Value = 39666.668
$ws.Range("M34").Value = -3147
$ws.Range("N34").Value = -40072.668
$ws.Range("H36").Value = 25140
$ws.Range("I36").Value = 3350
$ws.Range("J36").Value = 39666.668
$ws.Range("K36").Value = 3350
$ws.Range("L36").Value = 39666.668
$ws.Range("M36").Value = -2635
$ws.Range("N36").Value = -41096.668
$ws.Range("H38").Value = 475.25
$ws.Range("I38").Value = 117.5625
$ws.Range("J38").Value = 832.9375
$ws.Range("K38").Value = 352.6875
$ws.Range("L38").Value = 2498.8125
$ws.Range("M38").Value = 19.3125
$ws.Range("N38").Value = -3242.8125
$ws.Range("H39").Value = 254
$ws.Range("I39").Value = 53.857143
$ws.Range("J39").Value = 381.36365
$ws.Range("K39").Value = 161.571429
$ws.Range("L39").Value = 1144.09095
$ws.Range("M39").Value = 134.428571
$ws.Range("N39").Value = -1736.09095
$ws.Range("H40").Value = 11780
$ws.Range("I40").Value = 17916.666
$ws.Range("J40").Value = 2575
$ws.Range("K40").Value = 17916.666
$ws.Range("L40").Value = 2575
$ws.Range("M40").Value = -17741.666
$ws.Range("N40").Value = -2925
$ws.Range("H41").Value = 20478
$ws.Range("I41").Value = 50195
$ws.Range("K41").Value = 50195
$ws.Range("M41").Value = -49755
$ws.Range("H42").Value = 203.47368
$ws.Range("I42").Value = 45
$ws.Range("J42").Value = 318.72726
$ws.Range("K42").Value = 135
$ws.Range("L42").Value = 956.18178
$ws.Range("M42").Value = 95
$ws.Range("N42").Value = -1416.18178
$ws.Range("H43").Value = 790.3
$ws.Range("I43").Value = 978.2727
$ws.Range("J43").Value = 560.55554
$ws.Range("K43").Value = 978.2727
$ws.Range("L43").Value = 560.55554
$ws.Range("M43").Value = -909.2727
$ws.Range("N43").Value = -698.55554
$ws.Range("H47").Value = 62500
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 62500
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 62500
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -64444
$ws.Range("H48").Value = 1699.25
$ws.Range("I48").Value = 1100
$ws.Range("J48").Value = 1784.8572
$ws.Range("K48").Value = 3300
$ws.Range("L48").Value = 5354.571599999999
$ws.Range("M48").Value = -3008
$ws.Range("N48").Value = -5938.571599999999
$ws.Range("H51").Value = 1758.8823
$ws.Range("I51").Value = 1607.7693
$ws.Range("J51").Value = 2250
$ws.Range("K51").Value = 1607.7693
$ws.Range("L51").Value = 2250
$ws.Range("M51").Value = -1123.7693
$ws.Range("N51").Value = -3218
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()
$ws.Range("H53").Value = 192.33333
$ws.Range("I53").Value = 210.8
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 210.8
$ws.Range("L53").Value = 100
$ws.Range("M53").Value = 426.2
$ws.Range("N53").Value = -1374
$ws.Range("H54").Value = 50000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 50000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 50000
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -50972
$ws.Range("H55").Value = 407.65518
$ws.Range("I55").Value = 211.8
$ws.Range("J55").Value = 510.73685
$ws.Range("K55").Value = 211.8
$ws.Range("L55").Value = 510.73685
$ws.Range("M55").Value = 2.199999999999989
$ws.Range("N55").Value = -938.73685
$ws.Range("H56").Value = 1699.25
$ws.Range("I56").Value = 1100
$ws.Range("J56").Value = 1784.8572
$ws.Range("K56").Value = 3300
$ws.Range("L56").Value = 5354.571599999999
$ws.Range("M56").Value = -2766
$ws.Range("N56").Value = -6422.571599999999
$ws.Range("H58").Value = 1619.4546
$ws.Range("I58").Value = 288.8
$ws.Range("J58").Value = 2728.3333
$ws.Range("K58").Value = 866.4000000000001
$ws.Range("L58").Value = 8184.999899999999
$ws.Range("M58").Value = -716.4000000000001
$ws.Range("N58").Value = -8484.999899999999
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71248
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -216240
$ws.Range("H111").Value = 1983.3334
$ws.Range("I111").Value = 1300
$ws.Range("J111").Value = 2325
$ws.Range("K111").Value = 3900
$ws.Range("L111").Value = 6975
$ws.Range("M111").Value = -833
$ws.Range("N111").Value = -13109
$ws.Range("H116").Value = 1627294.8
$ws.Range("I116").Value = 8930296
$ws.Range("J116").Value = 4405.5557
$ws.Range("K116").Value = 8930296
$ws.Range("L116").Value = 4405.5557
$ws.Range("M116").Value = -8926854
$ws.Range("N116").Value = -11289.5557
$ws.Range("H138").Value = 3188.2812
$ws.Range("I138").Value = 1069.8077
$ws.Range("J138").Value = 4637.763
$ws.Range("K138").Value = 3209.4231
$ws.Range("L138").Value = 13913.289
$ws.Range("M138").Value = 1930.5769
$ws.Range("N138").Value = -24193.289

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2162.375
$ws.Range("I2").Value = 1099.6666
$ws.Range("J2").Value = 2800
$ws.Range("K2").Value = 1099.6666
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = -986.6666
$ws.Range("N2").Value = -3026
$ws.Range("H61").Value = 1779.6316
$ws.Range("I61").Value = 1636.0588
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1636.0588
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1424.0588
$ws.Range("N61").Value = -3424
$ws.Range("H116").Value = 2162.375
$ws.Range("I116").Value = 1099.6666
$ws.Range("J116").Value = 2800
$ws.Range("K116").Value = 1099.6666
$ws.Range("L116").Value = 2800
$ws.Range("M116").Value = 1194.3334
$ws.Range("N116").Value = -7388
$ws.Range("H132").Value = 2029.5
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 2700
$ws.Range("M132").Value = -170
$ws.Range("H136").Value = 1779.6316
$ws.Range("I136").Value = 1636.0588
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 4908.1764
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -2358.1764
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2162.375
$ws.Range("I3").Value = 1099.6666
$ws.Range("J3").Value = 2800
$ws.Range("K3").Value = 1099.6666
$ws.Range("L3").Value = 2800
$ws.Range("M3").Value = -985.6666
$ws.Range("N3").Value = -3028
$ws.Range("H94").Value = 757.1429000000001
$ws.Range("I94").Value = 757.1429000000001
$ws.Range("K94").Value = 757.1429000000001
$ws.Range("M94").Value = -306.1429000000001
$ws.Range("H134").Value = 1942.9231
$ws.Range("I134").Value = 1724.875
$ws.Range("J134").Value = 2939.7144
$ws.Range("K134").Value = 5174.625
$ws.Range("L134").Value = 8819.143199999999
$ws.Range("M134").Value = -2639.625
$ws.Range("N134").Value = -13889.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 5762.875
$ws.Range("I35").Value = 3739.2856
$ws.Range("J35").Value = 19928
$ws.Range("K35").Value = 3739.2856
$ws.Range("L35").Value = 19928
$ws.Range("M35").Value = -3445.2856
$ws.Range("N35").Value = -20516
$ws.Range("H99").Value = 1840
$ws.Range("I99").Value = 800
$ws.Range("K99").Value = 800
$ws.Range("M99").Value = 698
$ws.Range("H105").Value = 1353.625
$ws.Range("I105").Value = 1455.3636
$ws.Range("J105").Value = 1129.8
$ws.Range("K105").Value = 1455.3636
$ws.Range("L105").Value = 1129.8
$ws.Range("M105").Value = 291.6364000000001
$ws.Range("N105").Value = -4623.8
$ws.Range("H107").Value = 919.9048
$ws.Range("I107").Value = 442.33334
$ws.Range("J107").Value = 2113.8333
$ws.Range("K107").Value = 442.33334
$ws.Range("L107").Value = 2113.8333
$ws.Range("M107").Value = 1477.66666
$ws.Range("N107").Value = -5953.8333
$ws.Range("H122").Value = 1194.5714
$ws.Range("I122").Value = 1182.6666
$ws.Range("J122").Value = 1210.4445
$ws.Range("K122").Value = 3547.9998
$ws.Range("L122").Value = 3631.3335
$ws.Range("M122").Value = -1097.9998
$ws.Range("N122").Value = -8531.333500000001
$ws.Range("H126").Value = 1840
$ws.Range("I126").Value = 800
$ws.Range("K126").Value = 2400
$ws.Range("M126").Value = 70

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 380.4762
$ws.Range("I107").Value = 404
$ws.Range("J107").Value = 371.06668
$ws.Range("K107").Value = 1212
$ws.Range("L107").Value = 1113.20004
$ws.Range("M107").Value = 708
$ws.Range("N107").Value = -4953.20004
$ws.Range("H110").Value = 6595
$ws.Range("I110").Value = 3960
$ws.Range("J110").Value = 9230
$ws.Range("K110").Value = 11880
$ws.Range("L110").Value = 27690
$ws.Range("M110").Value = -7790
$ws.Range("N110").Value = -35870
$ws.Range("H122").Value = 4408.222
$ws.Range("I122").Value = 540.9286
$ws.Range("J122").Value = 8573
$ws.Range("K122").Value = 4868.3574
$ws.Range("L122").Value = 77157
$ws.Range("M122").Value = -2418.3574
$ws.Range("N122").Value = -82057

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2689.7273
$ws.Range("I122").Value = 2741.111
$ws.Range("K122").Value = 8223.332999999999
$ws.Range("M122").Value = -5773.332999999999
$ws.Range("H126").Value = 2263.8333
$ws.Range("I126").Value = 2233.3333
$ws.Range("J126").Value = 2294.3333
$ws.Range("K126").Value = 6699.999899999999
$ws.Range("L126").Value = 6882.999899999999
$ws.Range("M126").Value = -4229.999899999999
$ws.Range("N126").Value = -11822.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 9584.77
$ws.Range("I93").Value = 25750.75
$ws.Range("J93").Value = 2399.889
$ws.Range("K93").Value = 25750.75
$ws.Range("L93").Value = 2399.889
$ws.Range("N93").Value = -4895.889
